$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1 & 2. Find the paragraph holding the manual page-break run
#        (<w:r><w:lastRenderedPageBreak/><w:br w:type="page"/></w:r>) and,
#        right after it (but still inside the same paragraph, before the
#        paragraph mark), type the new sentence "Tämän pelin idea on" as its
#        own run, then drop a (new) "_GoBack" bookmark immediately after it.
#        Word only ever keeps a single "_GoBack" bookmark in a document, so
#        creating the new one automatically relocates/replaces the one that
#        used to sit after "Miska Sainkangas" on the cover page.
# ---------------------------------------------------------------------------

$breakParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "\f") {
        $breakParagraph = $p
    }
}

if ($breakParagraph -ne $null) {
    $nextParagraph = $breakParagraph.Next()

    # Type the sentence (plus one throw-away guard character) into the
    # paragraph that immediately follows the page break, so it becomes a
    # fresh run of its own rather than merging into the break's run.
    $nextRange = $nextParagraph.Range
    $nextRange.InsertAfter("Tämän pelin idea onX")

    # Copy just the typed text (no paragraph mark) and paste it at the end
    # of the page-break paragraph, right before its paragraph mark - this
    # keeps it as an independent run and leaves the original break run (and
    # its lastRenderedPageBreak) untouched.
    $nextParagraph = $breakParagraph.Next()
    $typed = $d.Range($nextParagraph.Range.Start, $nextParagraph.Range.End - 1)
    $typed.Copy()

    $pasteAt = $d.Range($breakParagraph.Range.End - 1, $breakParagraph.Range.End - 1)
    $pasteAt.Paste()

    # Remove the scratch copy that is still sitting in the next paragraph,
    # restoring it to its original (empty) state.
    $nextParagraph = $breakParagraph.Next()
    $scratch = $d.Range($nextParagraph.Range.Start, $nextParagraph.Range.End - 1)
    $scratch.Delete()

    # Bookmark a safe (non paragraph-mark) position right before the guard
    # character "X", then delete the guard - the collapsed bookmark stays
    # anchored right after "on", which is exactly where it belongs.
    $bpEnd = $breakParagraph.Range.End
    $bmRange = $d.Range($bpEnd - 2, $bpEnd - 2)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $guard = $d.Range($breakParagraph.Range.End - 2, $breakParagraph.Range.End - 1)
    $guard.Delete()
}

# ---------------------------------------------------------------------------
# 3. The extra paragraph pushes the document to a new page, so the cached
#    "PAGE" field result shown in the header used on page 2 onward needs to
#    go from 1 to 2.
# ---------------------------------------------------------------------------

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("1", $false, $false, $false, $false, $false,
                                     $true, 1, $false, "2", 2)
        }
    }
}
